$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agile User Story")
$ws.Activate()

# Fill in the next batch of story numbers in column B (rows 6-15),
# continuing the existing sequence (B3=1 ... B5=3) with the new
# folders/examples added for this commit.
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 8
$ws.Range("B11").Value = 9
$ws.Range("B12").Value = 10
$ws.Range("B13").Value = 11
$ws.Range("B14").Value = 12
$ws.Range("B15").Value = 13

# Scroll/select so the view lands on D16, matching where the author
# was working when they saved.
$ws.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("D16").Select()
